$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 372.5
$ws.Range("I38").Value = 200.2
$ws.Range("J38").Value = 1234
$ws.Range("K38").Value = 600.5999999999999
$ws.Range("L38").Value = 3702
$ws.Range("M38").Value = -228.5999999999999
$ws.Range("N38").Value = -4446

$ws.Range("H39").Value = 1325
$ws.Range("I39").Value = 433.33334
$ws.Range("K39").Value = 1300.00002
$ws.Range("M39").Value = -1004.00002

$ws.Range("H76").Value = 250006750
$ws.Range("I76").Value = 333337340
$ws.Range("K76").Value = 333337340
$ws.Range("M76").Value = -333337025

$ws.Range("H79").Value = 250006750
$ws.Range("I79").Value = 333337340
$ws.Range("K79").Value = 333337340
$ws.Range("M79").Value = -333336248

$ws.Range("H116").Value = 3845654.5
$ws.Range("I116").Value = 5569404
$ws.Range("K116").Value = 5569404
$ws.Range("M116").Value = -5565962

$ws.Range("H127").Value = 7177.9414
$ws.Range("I127").Value = 7850.5713
$ws.Range("K127").Value = 23551.7139
$ws.Range("M127").Value = -18591.7139

$ws.Range("H129").Value = 19622666
$ws.Range("I129").Value = 3269.7144
$ws.Range("J129").Value = 33356242
$ws.Range("K129").Value = 9809.143199999999
$ws.Range("L129").Value = 100068726
$ws.Range("M129").Value = -4809.143199999999
$ws.Range("N129").Value = -100078726

$ws.Range("H132").Value = 295318.22
$ws.Range("I132").Value = 322700.3
$ws.Range("K132").Value = 968100.8999999999
$ws.Range("M132").Value = -965570.8999999999

$ws.Range("H135").Value = 6298.9
$ws.Range("I135").Value = 3798
$ws.Range("K135").Value = 34182
$ws.Range("M135").Value = -31647

$ws.Range("H138").Value = 4940.385
$ws.Range("I138").Value = 1097.75
$ws.Range("J138").Value = 5639.0454
$ws.Range("K138").Value = 3293.25
$ws.Range("L138").Value = 16917.1362
$ws.Range("M138").Value = 1846.75
$ws.Range("N138").Value = -27197.1362

$ws.Range("H141").Value = 3397.4443
$ws.Range("I141").Value = 3771
$ws.Range("K141").Value = 11313
$ws.Range("M141").Value = -6133

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4464.086
$ws.Range("I2").Value = 4309.4
$ws.Range("K2").Value = 4309.4
$ws.Range("M2").Value = -4196.4

$ws.Range("H116").Value = 4464.086
$ws.Range("I116").Value = 4309.4
$ws.Range("K116").Value = 4309.4
$ws.Range("M116").Value = -2015.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4464.086
$ws.Range("I3").Value = 4309.4
$ws.Range("K3").Value = 4309.4
$ws.Range("M3").Value = -4195.4

$ws.Range("H86").Value = 10821.286
$ws.Range("I86").Value = 8333.333000000001
$ws.Range("J86").Value = 12687.25
$ws.Range("K86").Value = 8333.333000000001
$ws.Range("L86").Value = 12687.25
$ws.Range("M86").Value = -7210.333000000001
$ws.Range("N86").Value = -14933.25

$ws.Range("H89").Value = 10821.286
$ws.Range("I89").Value = 8333.333000000001
$ws.Range("J89").Value = 12687.25
$ws.Range("K89").Value = 41666.665
$ws.Range("L89").Value = 63436.25
$ws.Range("M89").Value = -36050.665
$ws.Range("N89").Value = -74668.25

$ws.Range("H99").Value = 4842.16
$ws.Range("I99").Value = 2154.5
$ws.Range("J99").Value = 8262.817999999999
$ws.Range("K99").Value = 2154.5
$ws.Range("L99").Value = 8262.817999999999
$ws.Range("M99").Value = -656.5
$ws.Range("N99").Value = -11258.818

$ws.Range("H134").Value = 3972986.8
$ws.Range("I134").Value = 4331463
$ws.Range("K134").Value = 12994389
$ws.Range("M134").Value = -12991854

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10770.857
$ws.Range("I62").Value = 11382.667
$ws.Range("J62").Value = 10312
$ws.Range("K62").Value = 11382.667
$ws.Range("L62").Value = 10312
$ws.Range("M62").Value = -10758.667
$ws.Range("N62").Value = -11560

$ws.Range("H65").Value = 10770.857
$ws.Range("I65").Value = 11382.667
$ws.Range("J65").Value = 10312
$ws.Range("K65").Value = 56913.335
$ws.Range("L65").Value = 51560
$ws.Range("M65").Value = -53793.335
$ws.Range("N65").Value = -57800

$ws.Range("H74").Value = 22221.75
$ws.Range("I74").Value = 22221
$ws.Range("K74").Value = 22221
$ws.Range("M74").Value = -21347

$ws.Range("H77").Value = 22221.75
$ws.Range("I77").Value = 22221
$ws.Range("K77").Value = 66663
$ws.Range("M77").Value = -62295

$ws.Range("H88").Value = 8728.666999999999
$ws.Range("J88").Value = 8728.666999999999
$ws.Range("L88").Value = 8728.666999999999
$ws.Range("N88").Value = -9540.666999999999

$ws.Range("H91").Value = 8728.666999999999
$ws.Range("J91").Value = 8728.666999999999
$ws.Range("L91").Value = 8728.666999999999
$ws.Range("N91").Value = -11536.667

$ws.Range("H122").Value = 3484.5334
$ws.Range("I122").Value = 3591.2
$ws.Range("J122").Value = 3431.2
$ws.Range("K122").Value = 10773.6
$ws.Range("L122").Value = 10293.6
$ws.Range("M122").Value = -8323.599999999999
$ws.Range("N122").Value = -15193.6

$ws.Range("H132").Value = 4957.2974
$ws.Range("I132").Value = 4068.9644
$ws.Range("J132").Value = 7721
$ws.Range("K132").Value = 12206.8932
$ws.Range("L132").Value = 23163
$ws.Range("M132").Value = -9676.893199999999
$ws.Range("N132").Value = -28223

$ws.Range("H134").Value = 25007740
$ws.Range("J134").Value = 11770.048
$ws.Range("L134").Value = 35310.144
$ws.Range("N134").Value = -40380.144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 33335480
$ws.Range("I131").Value = 133333840
$ws.Range("J131").Value = 15153960
$ws.Range("K131").Value = 400001520
$ws.Range("L131").Value = 45461880
$ws.Range("M131").Value = -399996480
$ws.Range("N131").Value = -45471960

$ws.Range("H134").Value = 6571.5557
$ws.Range("I134").Value = 6571.5557
$ws.Range("K134").Value = 19714.6671
$ws.Range("M134").Value = -14644.6671

$ws.Range("H136").Value = 12822534
$ws.Range("I136").Value = 16668485
$ws.Range("J136").Value = 2699.6667
$ws.Range("K136").Value = 50005455
$ws.Range("L136").Value = 8099.000100000001
$ws.Range("M136").Value = -50000355
$ws.Range("N136").Value = -18299.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 58826464
$ws.Range("I132").Value = 66669340
$ws.Range("K132").Value = 200008020
$ws.Range("M132").Value = -200005490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2046.6061
$ws.Range("I55").Value = 906.2632
$ws.Range("J55").Value = 3594.2144
$ws.Range("K55").Value = 906.2632
$ws.Range("L55").Value = 3594.2144
$ws.Range("M55").Value = -733.2632
$ws.Range("N55").Value = -3940.2144

$ws.Range("H61").Value = 6086.2383
$ws.Range("I61").Value = 5115.0586
$ws.Range("K61").Value = 5115.0586
$ws.Range("M61").Value = -4913.0586

$ws.Range("H68").Value = 1976.6875
$ws.Range("I68").Value = 1844.7858
$ws.Range("J68").Value = 2900
$ws.Range("K68").Value = 1844.7858
$ws.Range("L68").Value = 2900
$ws.Range("M68").Value = -1095.7858
$ws.Range("N68").Value = -4398

$ws.Range("H71").Value = 1976.6875
$ws.Range("I71").Value = 1844.7858
$ws.Range("J71").Value = 2900
$ws.Range("K71").Value = 9223.929
$ws.Range("L71").Value = 14500
$ws.Range("M71").Value = -5479.929
$ws.Range("N71").Value = -21988

$ws.Range("H113").Value = 6086.2383
$ws.Range("I113").Value = 5115.0586
$ws.Range("K113").Value = 5115.0586
$ws.Range("M113").Value = -2945.0586

$ws.Range("H132").Value = 3742.25
$ws.Range("J132").Value = 3949.75
$ws.Range("L132").Value = 11849.25
$ws.Range("N132").Value = -16909.25

$ws.Range("H136").Value = 12822812
$ws.Range("I136").Value = 29412748
$ws.Range("J136").Value = 3315.4546
$ws.Range("K136").Value = 88238244
$ws.Range("L136").Value = 9946.363799999999
$ws.Range("M136").Value = -88235694
$ws.Range("N136").Value = -15046.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3510.8333
$ws.Range("I126").Value = 1920.2142
$ws.Range("J126").Value = 5737.7
$ws.Range("K126").Value = 5760.642599999999
$ws.Range("L126").Value = 17213.1
$ws.Range("M126").Value = -3290.642599999999
$ws.Range("N126").Value = -22153.1

$ws.Range("H136").Value = 18561676
$ws.Range("J136").Value = 18675.715
$ws.Range("L136").Value = 56027.145
$ws.Range("N136").Value = -61127.145
